$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "image_url" column header (new shared string, index 5)
$ws.Range("F1").Value = "image_url"

# Widen column E to match the template's custom column width
$ws.Columns.Item(5).ColumnWidth = 12.7265625

# Move/restore the active selection to the new last header cell
$ws.Range("F1").Select() | Out-Null
